$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to text format while we write values that otherwise
# would be auto-converted to numbers by Excel (e.g. "608.82" -> 608.82).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.418.06"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "3.698.14"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "608.82"
$ws.Range("E5").Value = "  +4.83%  "

$ws.Range("D6").Value = "193.90"
$ws.Range("E6").Value = "  +14.14%  "

$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  +1.76%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").Value = "0.721"
$ws.Range("E9").Value = "  +3.14%  "

$ws.Range("D10").Value = "59.97"
$ws.Range("E10").Value = "  +17.06%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "0.0000286"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "10.39"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").Value = "4.273.86"
$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").Value = "3.680.02"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").Value = "19.34"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").Value = "68.198.88"
$ws.Range("E20").Value = "  +1.25%  "

$ws.Range("D21").Value = "407.04"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").Value = "4.57"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").Value = "89.86"
$ws.Range("E23").Value = "  +2.91%  "

$ws.Range("D24").Value = "11.62"
$ws.Range("E24").Value = "  +7.93%  "

$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").Value = "12.97"
$ws.Range("E26").Value = "  +2.32%  "

$ws.Range("D27").Value = "6.03"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("D28").Value = "3.76"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "9.59"
$ws.Range("E29").Value = "  +2.21%  "

$ws.Range("D30").Value = "32.52"
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").Value = "7.68"
$ws.Range("E31").Value = "  +2.78%  "

$ws.Range("D32").Value = "47.67"
$ws.Range("E32").Value = "  +10.87%  "

$ws.Range("D33").Value = "12.69"
$ws.Range("E33").Value = "  +2.85%  "

$ws.Range("E34").Value = "  +5.16%  "

$ws.Range("D35").Value = "629.07"
$ws.Range("E35").Value = "  +6.42%  "

$ws.Range("D36").Value = "67.39"
$ws.Range("E36").Value = "  +4.77%  "

$ws.Range("E37").Value = "  -6.74%  "

$ws.Range("D38").Value = "0.410"
$ws.Range("E38").Value = "  +4.33%  "

$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("E41").Value = "  +3.41%  "

$ws.Range("D42").Value = "3.00"
$ws.Range("E42").Value = "  +1.47%  "

$ws.Range("D43").Value = "0.0442"
$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("D44").Value = "2.63"
$ws.Range("E44").Value = "  -4.83%  "

$ws.Range("D45").Value = "2.882.54"
$ws.Range("E45").Value = "  +4.78%  "

$ws.Range("E46").Value = "  +5.11%  "

$ws.Range("D47").Value = "9.19"
$ws.Range("E47").Value = "  +0.80%  "

$ws.Range("D48").Value = "145.84"
$ws.Range("E48").Value = "  +3.29%  "

# Rows 49-50: coin order swap (dogwifhat <-> WEMIXToken)
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.66"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.62"
$ws.Range("E50").Value = "  -6.57%  "

# Row 51
$ws.Range("E51").Value = "  -2.90%  "

# Restore default (General) formatting / style for column D so the underlying
# cell style matches the original workbook (no explicit style index).
$priceRange.ClearFormats()
